# Apply crypto price/volume updates per commit "Updated cryptos list on Thu Oct 12 05:52:15 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.865.80"
$ws.Range("D3").Value = "1.564.00"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'206.00"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("E6").Value = "  -1.28%  "
$ws.Range("D8").Value = "'21.81"
$ws.Range("E8").Value = "  -2.00%  "
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("D10").Value = "'0.0585"
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").Value = "1.785.64"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "1.564.48"
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("E15").Value = "  -0.66%  "
$ws.Range("D16").Value = "26.872.41"
$ws.Range("D17").Value = "'61.31"
$ws.Range("E17").Value = "  -2.46%  "
$ws.Range("D18").Value = "'215.25"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("E19").Value = "  +1.61%  "
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("D23").Value = "'9.19"
$ws.Range("E23").Value = "  -1.95%  "
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("D25").Value = "'153.48"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("E26").Value = "  +1.32%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E29").Value = "  -1.27%  "
$ws.Range("E30").Value = "  +0.72%  "
$ws.Range("E31").Value = "  -3.39%  "
$ws.Range("E32").Value = "  -0.28%  "
$ws.Range("D33").Value = "1.399.74"
$ws.Range("E33").Value = "  +0.88%  "
$ws.Range("D34").Value = "'2.93"
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("E35").Value = "  -1.09%  "
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("D37").Value = "'0.911"
$ws.Range("E37").Value = "  -4.02%  "
$ws.Range("E38").Value = "  -1.08%  "
$ws.Range("E39").Value = "  +1.93%  "
$ws.Range("D40").Value = "'0.812"
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("E42").Value = "  +1.31%  "
$ws.Range("E43").Value = "  +6.54%  "
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("D46").Value = "'63.50"
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("D47").Value = "1.698.90"
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("D50").Value = "0.0₇0974"
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("D51").Value = "'0.0950"
$ws.Range("E51").Value = "  +0.82%  "
